# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-name suffixes to "_FV2410" / "_FV2504"
# 2) Wrap the data range in a table (Table1)
# 3) Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headersFV2410 = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)
for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headersFV2410[$i]
}

$headersFV2504 = @(
  "Segmentname_FV2504",
  "Segmentgruppe_FV2504",
  "Segment_FV2504",
  "Datenelement_FV2504",
  "Segment ID_FV2504",
  "Code_FV2504",
  "Qualifier_FV2504",
  "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504",
  "Bedingung_FV2504"
)
for ($i = 0; $i -lt $headersFV2504.Length; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $headersFV2504[$i]
}

$tblRange = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $tblRange, $null, 1)
$tbl.Name = "Table1"

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
